$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New animation rows for the Monkey mesh (matching Sailor mesh update)
$data = @(
    @(6, "finch",     354, 384),
    @(7, "groggy",    385, 438),
    @(8, "knockBack", 439, 523)
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row, 3).Value = $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]
    $ws.Cells.Item($row, 5).Value = $entry[2]
    $ws.Cells.Item($row, 6).Value = $entry[3]

    # copy formatting from the row above so new rows look consistent
    $ws.Range("C$($row-1):F$($row-1)").Copy()
    $ws.Range("C$($row):F$($row)").PasteSpecial(-4122)

    $row++
}

$ws.Range("F14").Select()
